$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Part List Report")
Write-Host $ws.Name
